$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = "@"
$r.Value = '28.935.63'
$r = $ws.Range('E2')
$r.NumberFormat = "@"
$r.Value = '  -1.90%  '
$r = $ws.Range('D3')
$r.NumberFormat = "@"
$r.Value = '1.901.54'
$r = $ws.Range('E3')
$r.NumberFormat = "@"
$r.Value = '  -3.74%  '
$r = $ws.Range('E4')
$r.NumberFormat = "@"
$r.Value = '  -0.01%  '
$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '324.36'
$r = $ws.Range('E5')
$r.NumberFormat = "@"
$r.Value = '  -0.78%  '
$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '1.002'
$r = $ws.Range('E6')
$r.NumberFormat = "@"
$r.Value = '  -0.14%  '
$r = $ws.Range('D7')
$r.NumberFormat = "@"
$r.Value = '0.4588'
$r = $ws.Range('E7')
$r.NumberFormat = "@"
$r.Value = '  -1.60%  '
$r = $ws.Range('D8')
$r.NumberFormat = "@"
$r.Value = '0.3817'
$r = $ws.Range('E8')
$r.NumberFormat = "@"
$r.Value = '  -2.60%  '
$r = $ws.Range('D9')
$r.NumberFormat = "@"
$r.Value = '0.07698'
$r = $ws.Range('E9')
$r.NumberFormat = "@"
$r.Value = '  -3.26%  '
$r = $ws.Range('E10')
$r.NumberFormat = "@"
$r.Value = '  -1.60%  '
$r = $ws.Range('D11')
$r.NumberFormat = "@"
$r.Value = '22.03'
$r = $ws.Range('E11')
$r.NumberFormat = "@"
$r.Value = '  -3.33%  '
$r = $ws.Range('D12')
$r.NumberFormat = "@"
$r.Value = '1.923.09'
$r = $ws.Range('E12')
$r.NumberFormat = "@"
$r.Value = '  -1.46%  '
$r = $ws.Range('D13')
$r.NumberFormat = "@"
$r.Value = '6.932'
$r = $ws.Range('E13')
$r.NumberFormat = "@"
$r.Value = '  -3.56%  '
$r = $ws.Range('D14')
$r.NumberFormat = "@"
$r.Value = '5.650'
$r = $ws.Range('E14')
$r.NumberFormat = "@"
$r.Value = '  -3.36%  '
$r = $ws.Range('D15')
$r.NumberFormat = "@"
$r.Value = '0.07040'
$r = $ws.Range('E15')
$r.NumberFormat = "@"
$r.Value = '  -0.44%  '
$r = $ws.Range('E16')
$r.NumberFormat = "@"
$r.Value = '  -0.14%  '
$r = $ws.Range('D17')
$r.NumberFormat = "@"
$r.Value = '83.63'
$r = $ws.Range('E17')
$r.NumberFormat = "@"
$r.Value = '  -4.63%  '
$r = $ws.Range('D18')
$r.NumberFormat = "@"
$r.Value = '0.000009451'
$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '16.62'
$r = $ws.Range('E19')
$r.NumberFormat = "@"
$r.Value = '  -4.01%  '
$r = $ws.Range('D20')
$r.NumberFormat = "@"
$r.Value = '1.003'
$r = $ws.Range('E20')
$r.NumberFormat = "@"
$r.Value = '  -0.11%  '
$r = $ws.Range('D21')
$r.NumberFormat = "@"
$r.Value = '28.948.03'
$r = $ws.Range('E21')
$r.NumberFormat = "@"
$r.Value = '  -1.83%  '
$r = $ws.Range('D22')
$r.NumberFormat = "@"
$r.Value = '5.292'
$r = $ws.Range('E22')
$r.NumberFormat = "@"
$r.Value = '  -4.56%  '
$r = $ws.Range('E23')
$r.NumberFormat = "@"
$r.Value = '  -2.79%  '
$r = $ws.Range('B24')
$r.NumberFormat = "@"
$r.Value = 'Toncoin'
$r = $ws.Range('C24')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$r = $ws.Range('D24')
$r.NumberFormat = "@"
$r.Value = '2.095'
$r = $ws.Range('E24')
$r.NumberFormat = "@"
$r.Value = '  -0.63%  '
$r = $ws.Range('B25')
$r.NumberFormat = "@"
$r.Value = 'Monero'
$r = $ws.Range('C25')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$r = $ws.Range('D25')
$r.NumberFormat = "@"
$r.Value = '158.12'
$r = $ws.Range('E25')
$r.NumberFormat = "@"
$r.Value = '  -0.30%  '
$r = $ws.Range('B26')
$r.NumberFormat = "@"
$r.Value = 'EthereumClassic'
$r = $ws.Range('C26')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$r = $ws.Range('D26')
$r.NumberFormat = "@"
$r.Value = '19.03'
$r = $ws.Range('E26')
$r.NumberFormat = "@"
$r.Value = '  -2.45%  '
$r = $ws.Range('B27')
$r.NumberFormat = "@"
$r.Value = 'InternetComputer(DFINITY)'
$r = $ws.Range('C27')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$r = $ws.Range('D27')
$r.NumberFormat = "@"
$r.Value = '5.639'
$r = $ws.Range('E27')
$r.NumberFormat = "@"
$r.Value = '  -2.61%  '
$r = $ws.Range('B28')
$r.NumberFormat = "@"
$r.Value = 'BitcoinCash'
$r = $ws.Range('C28')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$r = $ws.Range('D28')
$r.NumberFormat = "@"
$r.Value = '117.40'
$r = $ws.Range('E28')
$r.NumberFormat = "@"
$r.Value = '  -1.93%  '
$r = $ws.Range('B29')
$r.NumberFormat = "@"
$r.Value = 'LidoDAOToken'
$r = $ws.Range('C29')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$r = $ws.Range('D29')
$r.NumberFormat = "@"
$r.Value = '1.837'
$r = $ws.Range('E29')
$r.NumberFormat = "@"
$r.Value = '  -3.80%  '
$r = $ws.Range('B30')
$r.NumberFormat = "@"
$r.Value = 'Stellar'
$r = $ws.Range('C30')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$r = $ws.Range('D30')
$r.NumberFormat = "@"
$r.Value = '0.09252'
$r = $ws.Range('E30')
$r.NumberFormat = "@"
$r.Value = '  -1.81%  '
$r = $ws.Range('B31')
$r.NumberFormat = "@"
$r.Value = 'ImmutableX'
$r = $ws.Range('C31')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$r = $ws.Range('D31')
$r.NumberFormat = "@"
$r.Value = '0.8628'
$r = $ws.Range('E31')
$r.NumberFormat = "@"
$r.Value = '  -3.52%  '
$r = $ws.Range('B32')
$r.NumberFormat = "@"
$r.Value = 'Filecoin'
$r = $ws.Range('C32')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$r = $ws.Range('D32')
$r.NumberFormat = "@"
$r.Value = '5.084'
$r = $ws.Range('E32')
$r.NumberFormat = "@"
$r.Value = '  -2.88%  '
$r = $ws.Range('B33')
$r.NumberFormat = "@"
$r.Value = 'ARBITRUM'
$r = $ws.Range('C33')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$r = $ws.Range('D33')
$r.NumberFormat = "@"
$r.Value = '1.240'
$r = $ws.Range('E33')
$r.NumberFormat = "@"
$r.Value = '  -6.32%  '
$r = $ws.Range('B34')
$r.NumberFormat = "@"
$r.Value = 'HuobiToken'
$r = $ws.Range('C34')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$r = $ws.Range('D34')
$r.NumberFormat = "@"
$r.Value = '2.996'
$r = $ws.Range('E34')
$r.NumberFormat = "@"
$r.Value = '  -5.61%  '
$r = $ws.Range('B35')
$r.NumberFormat = "@"
$r.Value = 'Hedera'
$r = $ws.Range('C35')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$r = $ws.Range('D35')
$r.NumberFormat = "@"
$r.Value = '0.05708'
$r = $ws.Range('E35')
$r.NumberFormat = "@"
$r.Value = '  -1.90%  '
$r = $ws.Range('B36')
$r.NumberFormat = "@"
$r.Value = 'TrustWalletToken'
$r = $ws.Range('C36')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$r = $ws.Range('D36')
$r.NumberFormat = "@"
$r.Value = '1.147'
$r = $ws.Range('E36')
$r.NumberFormat = "@"
$r.Value = '  -2.04%  '
$r = $ws.Range('B37')
$r.NumberFormat = "@"
$r.Value = 'Frax'
$r = $ws.Range('C37')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$r = $ws.Range('D37')
$r.NumberFormat = "@"
$r.Value = '1.002'
$r = $ws.Range('E37')
$r.NumberFormat = "@"
$r.Value = '  -0.07%  '
$r = $ws.Range('B38')
$r.NumberFormat = "@"
$r.Value = 'VeChain'
$r = $ws.Range('C38')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$r = $ws.Range('D38')
$r.NumberFormat = "@"
$r.Value = '0.02039'
$r = $ws.Range('E38')
$r.NumberFormat = "@"
$r.Value = '  -3.14%  '
$r = $ws.Range('B39')
$r.NumberFormat = "@"
$r.Value = 'TheSandbox'
$r = $ws.Range('C39')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$r = $ws.Range('D39')
$r.NumberFormat = "@"
$r.Value = '0.5483'
$r = $ws.Range('E39')
$r.NumberFormat = "@"
$r.Value = '  -4.19%  '
$r = $ws.Range('B40')
$r.NumberFormat = "@"
$r.Value = 'FraxShare'
$r = $ws.Range('C40')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$r = $ws.Range('D40')
$r.NumberFormat = "@"
$r.Value = '7.381'
$r = $ws.Range('E40')
$r.NumberFormat = "@"
$r.Value = '  -5.03%  '
$r = $ws.Range('B41')
$r.NumberFormat = "@"
$r.Value = 'Algorand'
$r = $ws.Range('C41')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$r = $ws.Range('D41')
$r.NumberFormat = "@"
$r.Value = '0.1751'
$r = $ws.Range('E41')
$r.NumberFormat = "@"
$r.Value = '  -2.71%  '
$r = $ws.Range('B42')
$r.NumberFormat = "@"
$r.Value = 'MXToken'
$r = $ws.Range('C42')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$r = $ws.Range('D42')
$r.NumberFormat = "@"
$r.Value = '2.769'
$r = $ws.Range('E42')
$r.NumberFormat = "@"
$r.Value = '  +0.68%  '
$r = $ws.Range('D43')
$r.NumberFormat = "@"
$r.Value = '9.258'
$r = $ws.Range('E43')
$r.NumberFormat = "@"
$r.Value = '  -4.11%  '
$r = $ws.Range('B44')
$r.NumberFormat = "@"
$r.Value = 'Decentraland'
$r = $ws.Range('C44')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$r = $ws.Range('D44')
$r.NumberFormat = "@"
$r.Value = '0.5159'
$r = $ws.Range('E44')
$r.NumberFormat = "@"
$r.Value = '  -3.68%  '
$r = $ws.Range('B45')
$r.NumberFormat = "@"
$r.Value = 'EnergySwap'
$r = $ws.Range('C45')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$r = $ws.Range('D45')
$r.NumberFormat = "@"
$r.Value = '11.26'
$r = $ws.Range('E45')
$r.NumberFormat = "@"
$r.Value = '  -3.95%  '
$r = $ws.Range('B46')
$r.NumberFormat = "@"
$r.Value = 'Cronos'
$r = $ws.Range('C46')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$r = $ws.Range('D46')
$r.NumberFormat = "@"
$r.Value = '0.06815'
$r = $ws.Range('E46')
$r.NumberFormat = "@"
$r.Value = '  -1.64%  '
$r = $ws.Range('B47')
$r.NumberFormat = "@"
$r.Value = 'RenderToken'
$r = $ws.Range('C47')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$r = $ws.Range('D47')
$r.NumberFormat = "@"
$r.Value = '2.061'
$r = $ws.Range('E47')
$r.NumberFormat = "@"
$r.Value = '  -6.37%  '
$r = $ws.Range('B48')
$r.NumberFormat = "@"
$r.Value = 'PEPE'
$r = $ws.Range('C48')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$r = $ws.Range('D48')
$r.NumberFormat = "@"
$r.Value = '0.000002591'
$r = $ws.Range('E48')
$r.NumberFormat = "@"
$r.Value = '  -18.42%  '
$r = $ws.Range('B49')
$r.NumberFormat = "@"
$r.Value = 'Quant'
$r = $ws.Range('C49')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$r = $ws.Range('D49')
$r.NumberFormat = "@"
$r.Value = '110.32'
$r = $ws.Range('E49')
$r.NumberFormat = "@"
$r.Value = '  -2.97%  '
$r = $ws.Range('E50')
$r.NumberFormat = "@"
$r.Value = '  -3.47%  '
$r = $ws.Range('B51')
$r.NumberFormat = "@"
$r.Value = 'PaxDollar'
$r = $ws.Range('C51')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$r = $ws.Range('D51')
$r.NumberFormat = "@"
$r.Value = '1.002'
$r = $ws.Range('E51')
$r.NumberFormat = "@"
$r.Value = '  -0.22%  '
